$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.657.56"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.62%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.326.87"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.65%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "581.07"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.76%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "174.54"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.45%  "
$ws.Range("E8").Value = "  +2.25%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "3.324.52"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.72%  "
$ws.Range("E10").Value = "  +6.63%  "
$ws.Range("E11").Value = "  +1.73%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "46.77"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +4.97%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000271"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.54%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "695.99"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +7.12%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.867.52"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.53%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "8.37"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.97%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "67.632.53"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.52%  "
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.357.18"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.38%  "
$ws.Range("B19").Value = "TRON"
$ws.Range("C19").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.119"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.54%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.49"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.14%  "
$ws.Range("E21").Value = "  +3.54%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.890"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.03%  "
$ws.Range("E23").Value = "  +4.51%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "16.82"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.23%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "101.26"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +5.19%  "
$ws.Range("E26").Value = "  +1.86%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.67"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.87%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.37"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.66%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "32.92"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.26%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.50"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.15%  "
$ws.Range("E31").Value = "  +2.68%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "568.66"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.14%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "10.97"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.28%  "
$ws.Range("E34").Value = "  +3.39%  "
$ws.Range("B35").Value = "Dai"
$ws.Range("C35").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.24%  "
$ws.Range("B36").Value = "Maker"
$ws.Range("C36").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.715.46"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.33%  "
$ws.Range("B37").Value = "OKB"
$ws.Range("C37").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "57.23"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.69%  "
$ws.Range("E38").Value = "  -2.75%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "35.16"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +9.84%  "
$ws.Range("E40").Value = "  +4.03%  "
$ws.Range("E41").Value = "  +5.74%  "
$ws.Range("E42").Value = "  +1.03%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.32"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.65%  "
$ws.Range("E44").Value = "  +3.50%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0₃0667"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.65%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0405"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.28%  "
$ws.Range("E47").Value = "  +4.66%  "
$ws.Range("E48").Value = "  +1.96%  "
$ws.Range("E49").Value = "  -0.20%  "
$ws.Range("E50").Value = "  +0.32%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "131.48"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.41%  "
